{"js": "// Locate the (single) table in the document body.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// --- 1) Delete rows \"3\"..\"14\" (the 12 empty trailing rows), keeping only\n//        the header row, row \"1\" and row \"2\". Delete from the bottom up so\n//        earlier indices stay valid while we iterate. ---\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (let i = rows.items.length - 1; i >= 3; i--) {\n  rows.items[i].delete();\n}\nawait context.sync();\n\n// --- 2) Resize the 3 columns (applies to every remaining row at once). ---\ntable.getCell(0, 0).columnWidth = 440 / 20; // dxa -> points\ntable.getCell(0, 1).columnWidth = 5769 / 20;\ntable.getCell(0, 2).columnWidth = 2853 / 20;\nawait context.sync();\n\n// --- 3) Re-split the runs of row \"2\"'s \"Improvement\" cell and flag \"gpio\"\n//        with proofErr spell-check markers, keeping the visible text the\n//        same (\"Possibly breakout more gpio pins \"). ---\nconst improvementCell = table.getCell(2, 1);\nconst firstPara = improvementCell.body.paragraphs.getFirst();\nconst paraRange = firstPara.getRange();\n\nconst newParaXml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t>Possibly b</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">reakout more </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>gpio</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> pins </w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nparaRange.insertOoxml(newParaXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# --- 1) Delete rows \"3\"..\"14\" (the 12 empty trailing rows), keeping only\n#        the header row, row \"1\" and row \"2\". Delete from the bottom up so\n#        earlier row indices stay valid while iterating. ---\nfor ($i = $table.Rows.Count; $i -ge 4; $i--) {\n    $table.Rows.Item($i).Delete()\n}\n\n# --- 2) Resize the 3 columns (Cell.Width sets the whole column's width). ---\n$table.Cell(1, 1).Width = 440 / 20\n$table.Cell(1, 2).Width = 5769 / 20\n$table.Cell(1, 3).Width = 2853 / 20\n\n# --- 3) Re-split the runs of row \"2\"'s \"Improvement\" cell and flag \"gpio\"\n#        with proofErr spell-check markers, keeping the visible text the\n#        same (\"Possibly breakout more gpio pins \"). ---\n$improvementCell = $table.Cell(3, 2)\n$newParaXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:r><w:t>Possibly b</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">reakout more </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>gpio</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> pins </w:t></w:r>' +\n    '</w:p>'\n[void]$improvementCell.Range.InsertXML($newParaXml)\n"}
